$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Data"
$ws.Range("B1").Value = "Attivita’ principale"
$ws.Range("C1").Value = "Ore"
$ws.Range("D1").Value = "Commenti"

$ws.Range("A2").Value = "18/06 – 25/06 "
$ws.Range("B2").Value = "Consegna materiale, studio su Python, Raspberry e GPIO"
$ws.Range("C2").Value = 40

$ws.Range("A3").Value = "26/06 – 30/06"
$ws.Range("B3").Value = "Test di comunicazione, malfunzionamento di GPIO poi risolto"
$ws.Range("C3").Value = 35
$ws.Range("D3").Value = "Questo errore e’ stato un grande time sink"

$ws.Range("A4").Value = "04/07 – 11/07"
$ws.Range("B4").Value = "Inizializzazione di repository, prime prove su file transfer"
$ws.Range("C4").Value = 25

$ws.Range("A5").Value = "12/07 – 18/07"
$ws.Range("B5").Value = "Studio di checksum, stesura di documentazione, file transfer con checksum"
$ws.Range("C5").Value = 25

$ws.Range("A6").Value = "23/07 – 30/07"
$ws.Range("B6").Value = "Ricerche su network monitoring, studio di zabbix, e cacti"
$ws.Range("C6").Value = 20

$ws.Range("A7").Value = "03/08 – 07/08"
$ws.Range("B7").Value = "Ulteriore studio su zabbix, prime prove di export di dati"
$ws.Range("C7").Value = 20

$ws.Range("A8").Value = "09/08 – 14/08"
$ws.Range("B8").Value = "Studio su zabbix e simulazione, inizialmente possibilita’ di usare zax"
$ws.Range("C8").Value = 25

$ws.Range("A9").Value = "23/08 – 30/08"
$ws.Range("B9").Value = "Studio di come Zax funzionerebbe"
$ws.Range("C9").Value = 50
$ws.Range("D9").Value = "Potevo ottimizzare meglio il tempo, altro grande time sink per una priorita’ sbagliata."

$ws.Range("A10").Value = "04/09 – 07/09"
$ws.Range("B10").Value = "Passaggio ad oggetti trap"
$ws.Range("C10").Value = 25

$ws.Range("A11").Value = "10/09 – 19/09"
$ws.Range("B11").Value = "Lavoro su script di export + prova a script completo di configurazione"
$ws.Range("C11").Value = 70
$ws.Range("D11").Value = "Avendo studiato questi meglio il completamento di questa parte e’ stato molto piu’ veloce"

$ws.Range("A12").Value = "19/09 – 26/09"
$ws.Range("B12").Value = "Integrazione con file sender"
$ws.Range("C12").Value = 25

$ws.Range("A13").Value = "27/09 – 1/10"
$ws.Range("B13").Value = "Ultimi test"
$ws.Range("C13").Value = 20

$ws.Range("B14").Value = "Ore totali:"

$ws.Range("C14").Formula = "=SUM(C2:C13)"

$ws.Columns.Item(2).ColumnWidth = 80.5

$ws.Range("C12").Select() | Out-Null
